$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header description text (row 2) -----------------------------
# Column L is "modmass2" -> description now talks about peptide2.
# Column I is "modmass1" -> description now talks about peptide1.
# (L2 is set first so the two freshly-introduced shared strings land in the
#  same order as the target workbook: ...peptide2... before ...peptide1...)
$ws.Range("L2").Value2 = "mass of a modification within peptide2 (;-delimited if multiple)"
$ws.Range("I2").Value2 = "mass of a modification within peptide1 (;-delimited if multiple)"

# --- Update the worksheet's view / current selection ---------------------
$ws.Range("I3").Select()
